# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the newly scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 3;  Value = 2138 },
    @{ Row = 5;  Value = 11156 },
    @{ Row = 8;  Value = 307 },
    @{ Row = 10; Value = 11062 },
    @{ Row = 11; Value = 446 },
    @{ Row = 12; Value = 1139 },
    @{ Row = 14; Value = 1722 },
    @{ Row = 15; Value = 5548 },
    @{ Row = 17; Value = 3435 }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$updates4 = @(
    @{ Row = 3;  Value = 2138 },
    @{ Row = 7;  Value = 11156 },
    @{ Row = 10; Value = 307 },
    @{ Row = 12; Value = 11062 },
    @{ Row = 13; Value = 446 },
    @{ Row = 14; Value = 1139 },
    @{ Row = 16; Value = 1722 },
    @{ Row = 17; Value = 5548 },
    @{ Row = 19; Value = 3435 }
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates4) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
